# Add a new "2020" data column (Q) to the SDG 1.5.1 indicator sheet,
# mirroring the existing formatting of column P (the "2019" column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column P's formatting (and values, as a starting point) into column Q
# for every row of the table (row 3 = top border row, row 4 = year headers,
# rows 5-34 = data rows). This keeps fonts/borders/number formats consistent
# with the rest of the table.
$ws.Range("P3:P34").Copy($ws.Range("Q3:Q34"))

# Year header
$ws.Cells.Item(4, 17).Value = 2020

# Data values for 2020 ("-" marks a missing/not-available data point, same
# convention used throughout the rest of the table)
$ws.Cells.Item(5, 17).Value = 51
$ws.Cells.Item(6, 17).Value = 29
$ws.Cells.Item(7, 17).Value = 22
$ws.Cells.Item(8, 17).Value = 5
$ws.Cells.Item(9, 17).Value = 3
$ws.Cells.Item(10, 17).Value = 2
$ws.Cells.Item(11, 17).Value = 15
$ws.Cells.Item(12, 17).Value = 9
$ws.Cells.Item(13, 17).Value = 5
$ws.Cells.Item(14, 17).Value = "-"
$ws.Cells.Item(15, 17).Value = "-"
$ws.Cells.Item(16, 17).Value = "-"
$ws.Cells.Item(17, 17).Value = "-"
$ws.Cells.Item(18, 17).Value = "-"
$ws.Cells.Item(19, 17).Value = "-"
$ws.Cells.Item(20, 17).Value = 7
$ws.Cells.Item(21, 17).Value = 7
$ws.Cells.Item(22, 17).Value = "-"
$ws.Cells.Item(23, 17).Value = "-"
$ws.Cells.Item(24, 17).Value = "-"
$ws.Cells.Item(25, 17).Value = "-"
$ws.Cells.Item(26, 17).Value = 24
$ws.Cells.Item(27, 17).Value = 10
$ws.Cells.Item(28, 17).Value = 14
$ws.Cells.Item(29, 17).Value = "-"
$ws.Cells.Item(30, 17).Value = "-"
$ws.Cells.Item(31, 17).Value = "-"
$ws.Cells.Item(32, 17).Value = "-"
$ws.Cells.Item(33, 17).Value = "-"
$ws.Cells.Item(34, 17).Value = "-"

# Keep the sheet's dimension / selection in sync with the newly used range.
$ws.Range("Q35").Select()
